# "Update countries & provincias Spain"
# Refresh the COVID case table on sheet "Pais":
#  - bump the "Datos actualizados" timestamp in A1 from 09:04 to 10:04
#  - insert fresh data for Rusia/Dinamarca/Eslovaquia/Gabon/Belice/Curazao,
#    which pushes the neighboring countries (and their existing figures)
#    down one row in the case-count ranking
#  - update case/recovered/death counts for the rows that shifted as a
#    result

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = 'Datos actualizados a 7 de Mayo de 2020 a las 10:04'
# Row 8
$ws.Range("A8").Value = 'Rusia'
$ws.Range("B8").Value = 177160
$ws.Range("C8").Value = 11231
$ws.Range("D8").Value = 23803
$ws.Range("E8").Value = 151732
$ws.Range("F8").Value = 2300
$ws.Range("G8").Value = 88
$ws.Range("H8").Value = 1625
# Row 9
$ws.Range("A9").Value = 'Francia'
$ws.Range("B9").Value = 174191
$ws.Range("D9").Value = 53972
$ws.Range("E9").Value = 94410
$ws.Range("F9").Value = 3147
$ws.Range("H9").Value = 25809
# Row 10
$ws.Range("A10").Value = 'Alemania'
$ws.Range("B10").Value = 168162
$ws.Range("D10").Value = 139900
$ws.Range("E10").Value = 20987
$ws.Range("F10").Value = 1884
$ws.Range("H10").Value = 7275
# Row 29
$ws.Range("B29").Value = 20939
$ws.Range("C29").Value = 741
$ws.Range("E29").Value = 19285
# Row 36
$ws.Range("D36").Value = 4862
$ws.Range("E36").Value = 9145
# Row 37
$ws.Range("E37").Value = 7451
$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 868
# Row 42
$ws.Range("A42").Value = 'Dinamarca'
$ws.Range("B42").Value = 10083
$ws.Range("C42").Value = 145
$ws.Range("D42").Value = 7493
$ws.Range("E42").Value = 2084
$ws.Range("F42").Value = 46
$ws.Range("H42").Value = 506
# Row 43
$ws.Range("A43").Value = 'Filipinas'
$ws.Range("B43").Value = 10004
$ws.Range("D43").Value = 1506
$ws.Range("E43").Value = 7840
$ws.Range("F43").Value = 31
$ws.Range("H43").Value = 658
# Row 48
$ws.Range("B48").Value = 7979
$ws.Range("C48").Value = 5
$ws.Range("D48").Value = 4214
$ws.Range("E48").Value = 3502
$ws.Range("F48").Value = 52
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 263
# Row 68
$ws.Range("B68").Value = 2958
$ws.Range("C68").Value = 55
$ws.Range("D68").Value = 980
$ws.Range("E68").Value = 1965
# Row 88
$ws.Range("A88").Value = 'Eslovaquia'
$ws.Range("B88").Value = 1445
$ws.Range("C88").Value = 16
$ws.Range("D88").Value = 806
$ws.Range("E88").Value = 613
$ws.Range("F88").Value = 4
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 26
# Row 89
$ws.Range("A89").Value = 'Senegal'
$ws.Range("B89").Value = 1433
$ws.Range("D89").Value = 493
$ws.Range("E89").Value = 928
$ws.Range("F89").Value = 6
$ws.Range("H89").Value = 12
# Row 90
$ws.Range("B90").Value = 1433
$ws.Range("C90").Value = 5
$ws.Range("D90").Value = 739
$ws.Range("E90").Value = 645
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 49
# Row 106
$ws.Range("D106").Value = 213
$ws.Range("E106").Value = 512
# Row 125
$ws.Range("A125").Value = 'Gabon'
$ws.Range("B125").Value = 439
$ws.Range("C125").Value = 42
$ws.Range("D125").Value = 99
$ws.Range("E125").Value = 332
$ws.Range("F125").Value = 1
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = 8
# Row 126
$ws.Range("A126").Value = 'Reunion'
$ws.Range("B126").Value = 425
$ws.Range("D126").Value = 300
$ws.Range("E126").Value = 125
$ws.Range("F126").Value = 2
$ws.Range("H126").Value = 0
# Row 132
$ws.Range("D132").Value = 265
$ws.Range("E132").Value = 51
# Row 191
$ws.Range("A191").Value = 'Belice'
$ws.Range("D191").Value = 16
$ws.Range("H191").Value = 2
# Row 192
$ws.Range("A192").Value = 'Nueva Caledonia'
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0
# Row 198
$ws.Range("A198").Value = 'Curazao'
$ws.Range("D198").Value = 13
$ws.Range("H198").Value = 1
# Row 199
$ws.Range("A199").Value = 'Dominica'
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 0

"done"